$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.564.20"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'1.593.11"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'210.27"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'0.508"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "'0.245"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'19.51"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "'1.605.54"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'64.35"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "'26.564.18"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'0.0" + [char]8323 + "0736"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'207.65"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").Value = "'8.91"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'145.01"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("E34").Value = "  +13.51%  "
$ws.Range("D35").Value = "'1.279.90"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'0.599"
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'0.817"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "'5.42"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "'62.44"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'1.729.58"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "'89.21"
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'0.102"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.50"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.02%  "
